# Fruta / hortaliza, semanal
# Update the weekly data rows (2-12) with new Fecha (D), Volumen (J),
# Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M)
# and Precio $/Kg (P) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44557; J = 400; K = 13000; L = 14000; M = 13500; P = 750 }
    3  = @{ D = 44998; J = 320; K = 17000; L = 18000; M = 17500; P = 972 }
    4  = @{ D = 45068; J = 400; K = 16000; L = 17000; M = 16500; P = 917 }
    5  = @{ D = 44547; J = 200; K = 13000; L = 14000; M = 13500; P = 750 }
    6  = @{ D = 44960; J = 400; K = 19500; L = 20000; M = 19750; P = 1097 }
    7  = @{ D = 44977; J = 400; K = 16500; L = 17000; M = 16750; P = 931 }
    8  = @{ D = 44984; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
    9  = @{ D = 45005; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
    10 = @{ D = 44957; J = 400; K = 21000; L = 22000; M = 21500; P = 1194 }
    11 = @{ D = 44964; J = 300; K = 20000; L = 21000; M = 20500; P = 1139 }
    12 = @{ D = 44568; J = 500; K = 15000; L = 16000; M = 15500; P = 861 }
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("D$row").Value = $values.D
    $ws.Range("J$row").Value = $values.J
    $ws.Range("K$row").Value = $values.K
    $ws.Range("L$row").Value = $values.L
    $ws.Range("M$row").Value = $values.M
    $ws.Range("P$row").Value = $values.P
}
